$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3350799.8
$ws.Range("J17").Value = 4020699.5
$ws.Range("L17").Value = 12062098.5
$ws.Range("N17").Value = -12062434.5
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = ""
$ws.Range("H62").Value = 6965
$ws.Range("I62").Value = 5870.8335
$ws.Range("K62").Value = 5870.8335
$ws.Range("M62").Value = -5246.8335
$ws.Range("H65").Value = 6965
$ws.Range("I65").Value = 5870.8335
$ws.Range("K65").Value = 29354.1675
$ws.Range("M65").Value = -26234.1675
$ws.Range("H137").Value = 3014.0571
$ws.Range("I137").Value = 2258.5862
$ws.Range("J137").Value = 6665.5
$ws.Range("K137").Value = 6775.758600000001
$ws.Range("L137").Value = 19996.5
$ws.Range("M137").Value = -4225.758600000001
$ws.Range("N137").Value = -25096.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1458.6774
$ws.Range("I2").Value = 832.6667
$ws.Range("K2").Value = 832.6667
$ws.Range("M2").Value = -719.6667
$ws.Range("H61").Value = 2504.3809
$ws.Range("I61").Value = 2413.1333
$ws.Range("J61").Value = 2732.5
$ws.Range("K61").Value = 2413.1333
$ws.Range("L61").Value = 2732.5
$ws.Range("M61").Value = -2201.1333
$ws.Range("N61").Value = -3156.5
$ws.Range("H92").Value = 66363.336
$ws.Range("J92").Value = 66363.336
$ws.Range("L92").Value = 66363.336
$ws.Range("N92").Value = -71355.336
$ws.Range("H116").Value = 1458.6774
$ws.Range("I116").Value = 832.6667
$ws.Range("K116").Value = 832.6667
$ws.Range("M116").Value = 1461.3333
$ws.Range("H122").Value = 1955.75
$ws.Range("I122").Value = 1609.3235
$ws.Range("J122").Value = 3918.8333
$ws.Range("K122").Value = 4827.970499999999
$ws.Range("L122").Value = 11756.4999
$ws.Range("M122").Value = -2377.970499999999
$ws.Range("N122").Value = -16656.4999
$ws.Range("H136").Value = 2504.3809
$ws.Range("I136").Value = 2413.1333
$ws.Range("J136").Value = 2732.5
$ws.Range("K136").Value = 7239.3999
$ws.Range("L136").Value = 8197.5
$ws.Range("M136").Value = -4689.3999
$ws.Range("N136").Value = -13297.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1458.6774
$ws.Range("I3").Value = 832.6667
$ws.Range("K3").Value = 832.6667
$ws.Range("M3").Value = -718.6667
$ws.Range("H20").Value = 3405.8667
$ws.Range("I20").Value = 3311
$ws.Range("J20").Value = 3488.875
$ws.Range("K20").Value = 3311
$ws.Range("L20").Value = 3488.875
$ws.Range("M20").Value = -3064
$ws.Range("N20").Value = -3982.875
$ws.Range("H60").Value = 106945.336
$ws.Range("J60").Value = 106945.336
$ws.Range("L60").Value = 106945.336
$ws.Range("N60").Value = -108143.336
$ws.Range("H86").Value = 1701.95
$ws.Range("I86").Value = 2521.8572
$ws.Range("J86").Value = 1260.4615
$ws.Range("K86").Value = 2521.8572
$ws.Range("L86").Value = 1260.4615
$ws.Range("M86").Value = -1398.8572
$ws.Range("N86").Value = -3506.4615
$ws.Range("H89").Value = 1701.95
$ws.Range("I89").Value = 2521.8572
$ws.Range("J89").Value = 1260.4615
$ws.Range("K89").Value = 12609.286
$ws.Range("L89").Value = 6302.307499999999
$ws.Range("M89").Value = -6993.286
$ws.Range("N89").Value = -17534.3075
$ws.Range("H94").Value = 1054.1818
$ws.Range("I94").Value = 512.3125
$ws.Range("K94").Value = 512.3125
$ws.Range("M94").Value = -61.3125
$ws.Range("H105").Value = 2552.7878
$ws.Range("I105").Value = 1943.4348
$ws.Range("K105").Value = 1943.4348
$ws.Range("M105").Value = -196.4348
$ws.Range("H107").Value = 1587.5238
$ws.Range("I107").Value = 1397.3334
$ws.Range("K107").Value = 1397.3334
$ws.Range("M107").Value = 522.6666
$ws.Range("H134").Value = 11566390
$ws.Range("I134").Value = 2552426
$ws.Range("J134").Value = 47622250
$ws.Range("K134").Value = 7657278
$ws.Range("L134").Value = 142866750
$ws.Range("M134").Value = -7654743
$ws.Range("N134").Value = -142871820
$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 566.5714
$ws.Range("I7").Value = 533.2308
$ws.Range("K7").Value = 533.2308
$ws.Range("M7").Value = -420.2308
$ws.Range("H31").Value = 4020.7188
$ws.Range("I31").Value = 2174.0625
$ws.Range("K31").Value = 2174.0625
$ws.Range("M31").Value = -1879.0625
$ws.Range("H34").Value = 4020.7188
$ws.Range("I34").Value = 2174.0625
$ws.Range("K34").Value = 2174.0625
$ws.Range("M34").Value = -1972.0625
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H105").Value = 1577.5625
$ws.Range("J105").Value = 1665.8
$ws.Range("L105").Value = 1665.8
$ws.Range("N105").Value = -5159.8
$ws.Range("H107").Value = 35990.93
$ws.Range("I107").Value = 48381.668
$ws.Range("K107").Value = 48381.668
$ws.Range("M107").Value = -46461.668
$ws.Range("H122").Value = 1670.8636
$ws.Range("I122").Value = 1377.3334
$ws.Range("K122").Value = 4132.0002
$ws.Range("M122").Value = -1682.0002
$ws.Range("H141").Value = 1047231.4
$ws.Range("I141").Value = 73325
$ws.Range("J141").Value = 1155443.2
$ws.Range("K141").Value = 73325
$ws.Range("L141").Value = 1155443.2
$ws.Range("M141").Value = -68145
$ws.Range("N141").Value = -1165803.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 2379
$ws.Range("I47").Value = 2068.5
$ws.Range("J47").Value = 3000
$ws.Range("K47").Value = 6205.5
$ws.Range("L47").Value = 9000
$ws.Range("M47").Value = -5774.5
$ws.Range("N47").Value = -9862
$ws.Range("H57").Value = 4857.2856
$ws.Range("I57").Value = 2332.6667
$ws.Range("K57").Value = 6998.000100000001
$ws.Range("M57").Value = -6439.000100000001
$ws.Range("H114").Value = 8677
$ws.Range("I114").Value = 2000
$ws.Range("J114").Value = 12015.5
$ws.Range("K114").Value = 6000
$ws.Range("L114").Value = 36046.5
$ws.Range("M114").Value = -2746
$ws.Range("N114").Value = -42554.5
$ws.Range("H121").Value = 7144784.5
$ws.Range("I121").Value = 12500486
$ws.Range("J121").Value = 3849.6667
$ws.Range("K121").Value = 37501458
$ws.Range("L121").Value = 11549.0001
$ws.Range("M121").Value = -37500148
$ws.Range("N121").Value = -14169.0001
$ws.Range("H129").Value = 1631.625
$ws.Range("I129").Value = 1009.6667
$ws.Range("J129").Value = 2004.8
$ws.Range("K129").Value = 3029.0001
$ws.Range("L129").Value = 6014.4
$ws.Range("M129").Value = 1970.9999
$ws.Range("N129").Value = -16014.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2879.8
$ws.Range("I122").Value = 3631.3333
$ws.Range("K122").Value = 10893.9999
$ws.Range("M122").Value = -8443.999899999999
$ws.Range("H132").Value = 1845.4375
$ws.Range("I132").Value = 1902.2
$ws.Range("J132").Value = 994
$ws.Range("K132").Value = 5706.6
$ws.Range("L132").Value = 2982
$ws.Range("M132").Value = -3176.6
$ws.Range("N132").Value = -8042
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16490.385
$ws.Range("I7").Value = 16273.857
$ws.Range("K7").Value = 16273.857
$ws.Range("M7").Value = -16161.857
$ws.Range("H22").Value = 1482.75
$ws.Range("I22").Value = 1340.4546
$ws.Range("J22").Value = 1603.1538
$ws.Range("K22").Value = 1340.4546
$ws.Range("L22").Value = 1603.1538
$ws.Range("M22").Value = -1045.4546
$ws.Range("N22").Value = -2193.1538
$ws.Range("H27").Value = 1482.75
$ws.Range("I27").Value = 1340.4546
$ws.Range("J27").Value = 1603.1538
$ws.Range("K27").Value = 1340.4546
$ws.Range("L27").Value = 1603.1538
$ws.Range("M27").Value = -1233.4546
$ws.Range("N27").Value = -1817.1538
$ws.Range("H55").Value = 244.17647
$ws.Range("I55").Value = 165.85
$ws.Range("K55").Value = 165.85
$ws.Range("M55").Value = 7.150000000000006
$ws.Range("H68").Value = 3850.1
$ws.Range("I68").Value = 2920.2
$ws.Range("J68").Value = 4780
$ws.Range("K68").Value = 2920.2
$ws.Range("L68").Value = 4780
$ws.Range("M68").Value = -2171.2
$ws.Range("N68").Value = -6278
$ws.Range("H71").Value = 3850.1
$ws.Range("I71").Value = 2920.2
$ws.Range("J71").Value = 4780
$ws.Range("K71").Value = 14601
$ws.Range("L71").Value = 23900
$ws.Range("M71").Value = -10857
$ws.Range("N71").Value = -31388
$ws.Range("H93").Value = 1097.2941
$ws.Range("I93").Value = 831.0833
$ws.Range("J93").Value = 1736.2
$ws.Range("K93").Value = 831.0833
$ws.Range("L93").Value = 1736.2
$ws.Range("M93").Value = 416.9167
$ws.Range("N93").Value = -4232.2
$ws.Range("H126").Value = 16490.385
$ws.Range("I126").Value = 16273.857
$ws.Range("K126").Value = 48821.571
$ws.Range("M126").Value = -46351.571
$ws.Range("H132").Value = 3772
$ws.Range("I132").Value = 2978.7778
$ws.Range("J132").Value = 5199.8
$ws.Range("K132").Value = 8936.3334
$ws.Range("L132").Value = 15599.4
$ws.Range("M132").Value = -6406.3334
$ws.Range("N132").Value = -20659.4
$ws.Range("H136").Value = 4471.6772
$ws.Range("I136").Value = 3877.1667
$ws.Range("J136").Value = 5294.846
$ws.Range("K136").Value = 11631.5001
$ws.Range("L136").Value = 15884.538
$ws.Range("M136").Value = -9081.500100000001
$ws.Range("N136").Value = -20984.538
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 433.33334
$ws.Range("I113").Value = 516.3333
$ws.Range("J113").Value = 267.33334
$ws.Range("K113").Value = 1548.9999
$ws.Range("L113").Value = 802.0000200000001
$ws.Range("M113").Value = 621.0001
$ws.Range("N113").Value = -5142.00002
$ws.Range("H136").Value = 19222.2
$ws.Range("I136").Value = 1473.1666
$ws.Range("K136").Value = 4419.4998
$ws.Range("M136").Value = -1869.4998
